$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45160

# Row 3
$ws.Range("D3").Value = 45229
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 2000
$ws.Range("P3").Value = 667

# Row 4
$ws.Range("D4").Value = 45203
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 2500
$ws.Range("M4").Value = 2500
$ws.Range("P4").Value = 833

# Row 5
$ws.Range("D5").Value = 45203
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1500
$ws.Range("P5").Value = 500

# Row 6
$ws.Range("D6").Value = 45205
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2250
$ws.Range("P6").Value = 750

# Row 7
$ws.Range("D7").Value = 44838
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1200
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = 1250
$ws.Range("P7").Value = 417

# Row 8
$ws.Range("D8").Value = 44838
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("P8").Value = 333

# Row 9
$ws.Range("D9").Value = 45148
$ws.Range("J9").Value = 80

# Row 10
$ws.Range("D10").Value = 45148
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 2000
$ws.Range("P10").Value = 667

# Row 11
$ws.Range("D11").Value = 45145
$ws.Range("J11").Value = 60

# Row 12
$ws.Range("D12").Value = 45145
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 80

# Row 13
$ws.Range("D13").Value = 45135
$ws.Range("J13").Value = 70

# Row 14
$ws.Range("D14").Value = 45161
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("P14").Value = 833

# Row 15
$ws.Range("D15").Value = 45191
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2500
$ws.Range("P15").Value = 833

# Row 16
$ws.Range("D16").Value = 45175
$ws.Range("I16").Value = "Primera"
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("P16").Value = 833

# Row 17
$ws.Range("D17").Value = 45215
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 2000
$ws.Range("P17").Value = 667

# Row 18
$ws.Range("D18").Value = 45163

# Row 19
$ws.Range("D19").Value = 45134
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 2500
$ws.Range("M19").Value = 2500
$ws.Range("P19").Value = 833

# Row 20
$ws.Range("D20").Value = 45133
$ws.Range("J20").Value = 80

# Row 21
$ws.Range("D21").Value = 45146
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("P21").Value = 833

# Row 22
$ws.Range("D22").Value = 45146
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 2000
$ws.Range("P22").Value = 667

# Row 23
$ws.Range("D23").Value = 45225
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = 2000
$ws.Range("P23").Value = 667

# Row 24
$ws.Range("D24").Value = 45149
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 2500
$ws.Range("M24").Value = 2500
$ws.Range("P24").Value = 833

# Row 25
$ws.Range("D25").Value = 45149
$ws.Range("I25").Value = "Segunda"
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 2000
$ws.Range("P25").Value = 667

# Row 26
$ws.Range("D26").Value = 44846
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 1200
$ws.Range("L26").Value = 1300
$ws.Range("M26").Value = 1250
$ws.Range("P26").Value = 417

# Row 27
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 1000
$ws.Range("P27").Value = 333

# Row 28
$ws.Range("D28").Value = 45217
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 2000
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2250
$ws.Range("P28").Value = 750

# Row 29
$ws.Range("D29").Value = 45176

# Row 30
$ws.Range("D30").Value = 44832
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 1200
$ws.Range("L30").Value = 1300
$ws.Range("M30").Value = 1250
$ws.Range("P30").Value = 417

# Row 31
$ws.Range("D31").Value = 44832
$ws.Range("J31").Value = 150
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = 1000
$ws.Range("P31").Value = 333

# Row 32
$ws.Range("D32").Value = 45195
$ws.Range("J32").Value = 100

# Row 33
$ws.Range("D33").Value = 45219
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 2000
$ws.Range("M33").Value = 2250
$ws.Range("P33").Value = 750

# Row 34
$ws.Range("D34").Value = 45166
$ws.Range("J34").Value = 120

# Row 35
$ws.Range("D35").Value = 45173
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 2500
$ws.Range("L35").Value = 2500
$ws.Range("M35").Value = 2500
$ws.Range("P35").Value = 833
